# Update Efemp1-Egfr NATMI results per Dr Hou advice: rescaled stats and added "M1" target cluster rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowVals = @("ECs", "Efemp1", "Egfr", "ECs", 3.0, 1.0, 1.51350766666666692, 4.54052300000000031, 0.0094290855215386, 0.00947765464651767, 3.0, 1.0, 1.21105733333333299, 3.63317200000000007, 0.0147461456544675, 0.01598314554371009, 1.83294455877288898, 16.49650102895600057, 0.00013904266848904, 0.00015148273362831)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])2").Value = $rowVals[$i]
}

$rowVals = @("ECs", "Efemp1", "Egfr", "FAPs", 3.0, 1.0, 1.51350766666666692, 4.54052300000000031, 0.0094290855215386, 0.00947765464651767, 3.0, 1.0, 61.06015933333333123, 183.1804779999999937, 0.74348420874183185, 0.80585236279493078, 92.41501927888822365, 831.73517350999406972, 0.00701037618814019, 0.00763759039065062)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])3").Value = $rowVals[$i]
}

$rowVals = @("ECs", "Efemp1", "Egfr", "M1", 3.0, 1.0, 1.51350766666666692, 4.54052300000000031, 0.0094290855215386, 0.00947765464651767, 1.0, 0.33333333333333331, 0.13227800000000001, 0.39683400000000002, 0.00161065095862375, 0.00174576253992177, 0.20020376713133331, 1.80183390418199996, 0.00001518696563421, 0.00001654573444821)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])4").Value = $rowVals[$i]
}

$rowVals = @("ECs", "Efemp1", "Egfr", "M2", 3.0, 1.0, 1.51350766666666692, 4.54052300000000031, 0.0094290855215386, 0.00947765464651767, 3.0, 1.0, 0.65512199999999998, 1.96536599999999995, 0.00797693401257583, 0.00864608964966683, 0.99153216960199997, 8.92378952641800005, 0.00007521519300425, 0.00008194465174237)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])5").Value = $rowVals[$i]
}

$rowVals = @("ECs", "Efemp1", "Egfr", "sCs", 3.0, 1.0, 1.51350766666666692, 4.54052300000000031, 0.0094290855215386, 0.00947765464651767, 2.0, 1.0, 19.06842599999999877, 38.13685199999999753, 0.23218206063250099, 0.16777263947177051, 28.86020894226599864, 173.16125365359599186, 0.00218926450627091, 0.00159009113604816)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])6").Value = $rowVals[$i]
}

$rowVals = @("FAPs", "Efemp1", "Egfr", "ECs", 3.0, 1.0, 156.53356166666671356, 469.60068499999999858, 0.97519713474375314, 0.98022036540683311, 3.0, 1.0, 1.21105733333333299, 3.63317200000000007, 0.0147461456544675, 0.01598314554371009, 189.57111776920220336, 1706.14005992282000079, 0.01438039899075075, 0.0156670047652061)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])7").Value = $rowVals[$i]
}

$rowVals = @("FAPs", "Efemp1", "Egfr", "FAPs", 3.0, 1.0, 156.53356166666671356, 469.60068499999999858, 0.97519713474375314, 0.98022036540683311, 3.0, 1.0, 61.06015933333333123, 183.1804779999999937, 0.74348420874183185, 0.80585236279493078, 9557.96421638082392747, 86021.67794742742262315, 0.72504367009226089, 0.78991289752280691)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])8").Value = $rowVals[$i]
}

$rowVals = @("FAPs", "Efemp1", "Egfr", "M1", 3.0, 1.0, 156.53356166666671356, 469.60068499999999858, 0.97519713474375314, 0.98022036540683311, 1.0, 0.33333333333333331, 0.13227800000000001, 0.39683400000000002, 0.00161065095862375, 0.00174576253992177, 20.70594647014333134, 186.35351823128999627, 0.00157070219992216, 0.00171123199479568)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])9").Value = $rowVals[$i]
}

$rowVals = @("FAPs", "Efemp1", "Egfr", "M2", 3.0, 1.0, 156.53356166666671356, 469.60068499999999858, 0.97519713474375314, 0.98022036540683311, 3.0, 1.0, 0.65512199999999998, 1.96536599999999995, 0.00797693401257583, 0.00864608964966683, 102.54857998619000625, 922.9372198757099568, 0.00777908319310394, 0.00847507315573666)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])10").Value = $rowVals[$i]
}

$rowVals = @("FAPs", "Efemp1", "Egfr", "sCs", 3.0, 1.0, 156.53356166666671356, 469.60068499999999858, 0.97519713474375314, 0.98022036540683311, 2.0, 1.0, 19.06842599999999877, 38.13685199999999753, 0.23218206063250099, 0.16777263947177051, 2984.84863715726896771, 17909.09182294362108223, 0.22642328026771541, 0.16445415796828769)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])11").Value = $rowVals[$i]
}

$rowVals = @("sCs", "Efemp1", "Egfr", "ECs", 2.0, 1.0, 2.46771900000000022, 4.93543800000000044, 0.01537377973470835, 0.01030197994664929, 3.0, 1.0, 1.21105733333333299, 3.63317200000000007, 0.0147461456544675, 0.01598314554371009, 2.98854919155599985, 17.93129514933600177, 0.00022670399522771, 0.00016465804487568)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])12").Value = $rowVals[$i]
}

$rowVals = @("sCs", "Efemp1", "Egfr", "FAPs", 2.0, 1.0, 2.46771900000000022, 4.93543800000000044, 0.01537377973470835, 0.01030197994664929, 3.0, 1.0, 61.06015933333333123, 183.1804779999999937, 0.74348420874183185, 0.80585236279493078, 150.67931532989399557, 904.07589197936385972, 0.01143016246143085, 0.00830187488147332)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])13").Value = $rowVals[$i]
}

$rowVals = @("sCs", "Efemp1", "Egfr", "M1", 2.0, 1.0, 2.46771900000000022, 4.93543800000000044, 0.01537377973470835, 0.01030197994664929, 1.0, 0.33333333333333331, 0.13227800000000001, 0.39683400000000002, 0.00161065095862375, 0.00174576253992177, 0.326424933882, 1.95854960329200001, 0.00002476179306738, 0.00001798481067789)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])14").Value = $rowVals[$i]
}

$rowVals = @("sCs", "Efemp1", "Egfr", "M2", 2.0, 1.0, 2.46771900000000022, 4.93543800000000044, 0.01537377973470835, 0.01030197994664929, 3.0, 1.0, 0.65512199999999998, 1.96536599999999995, 0.00797693401257583, 0.00864608964966683, 1.61665700671799994, 9.69994204030799878, 0.00012263562646764, 0.0000890718421878)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])15").Value = $rowVals[$i]
}

$rowVals = @("sCs", "Efemp1", "Egfr", "sCs", 2.0, 1.0, 2.46771900000000022, 4.93543800000000044, 0.01537377973470835, 0.01030197994664929, 2.0, 1.0, 19.06842599999999877, 38.13685199999999753, 0.23218206063250099, 0.16777263947177051, 47.05551714029399335, 188.22206856117600182, 0.00356951585851477, 0.0017283903674346)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])16").Value = $rowVals[$i]
}
